$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 with new email / password values
$ws.Range("C5").Value = "bmkbihhxjvucvxk@gmail.com"
$ws.Range("D5").Value = "rvjycHYANC5"

# Remove the Status value for row 5 (F5 no longer populated)
$ws.Range("F5").Value = $null

# Delete rows 6 through 9 entirely (shrinks used range to A1:F5)
$ws.Range("A6:F9").Delete()

# Update the active selection to match the recorded cursor position
$ws.Range("G16").Select()
